$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new test entry ---
# Set values first (so leading "-" text keeps its literal content), then copy
# number/text formatting from the matching column in an already-filled row
# (row 3) so the cells land on the exact same shared style as the rest of the
# table (date format, time format, wrapped/quote-prefixed note text).

$ws.Range("E8").Value = "- After 10 min, black Rpi lock to gps with time correct within 1ms and red Rpi lock to gps with time correct within 4ms`n- "
$ws.Range("E3").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$ws.Range("D8").Value = "- With internet from laptop (laptop is connected to Ethernet)`n- Lora transmitter is connected to laptop via usb uart CP2102 and transmit 26 packages each reset time`n- Each Rpi is connected with an LoRa reciver via an USB UART at /dev/ttyUSB0 and an ublox neo 7 GPS`n- Both Rpi ntp can synchronize with each other via LAN`n- Test indoor, room E6.1, 2 GPS with next by the window`n- 2 LoRa receiver is put side by side"
$ws.Range("D3").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("B8").Value = 43210
$ws.Range("B3").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C8").Value = 0.69444444444444453
$ws.Range("C3").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Rows.Item(8).RowHeight = 144

# --- Row 9: new test entry ---
$ws.Range("D9").Value = "- With internet from laptop (laptop is connected to phone's Ethernet)`n- Lora transmitter is connected to laptop via usb uart CP2102 and transmit 26 packages each reset time`n- Each Rpi is connected with an LoRa reciver via an USB UART at /dev/ttyUSB0 and an ublox neo 7 GPS`n- Both Rpi ntp can synchronize with each other via LAN`n- Test outdoor, roof of E building`n- 2 LoRa receiver is put side by side"
$ws.Range("D3").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$ws.Range("B9").Value = 43210
$ws.Range("B3").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C9").Value = 0.71597222222222223
$ws.Range("C3").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Rows.Item(9).RowHeight = 199.5

# --- Minor row-height touch-ups elsewhere in the table (re-flowed because of
# the new content above) ---
$ws.Rows.Item(3).RowHeight = 106.9
$ws.Rows.Item(4).RowHeight = 106.9
$ws.Rows.Item(5).RowHeight = 79.15
$ws.Rows.Item(6).RowHeight = 105
$ws.Rows.Item(7).RowHeight = 86.45
$ws.Rows.Item(10).RowHeight = 49.9
$ws.Rows.Item(11).RowHeight = 49.9
$ws.Rows.Item(12).RowHeight = 49.9

# --- View state: scrolled down to row 9, D10 selected ---
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D10").Select()
